# "Subindo algoritmo com reprodução"
# The genetic-algorithm parameters in column A (sheet "Planilha1") are
# updated: the population size configuration is bumped from 10 to 100,
# and a new configuration row is added for the number of generations.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A2: "Tamanho população: 10" -> "Tamanho população: 100"
$ws.Range("A2").Value = "Tamanho população: 100"

# A7 was empty; add the new "Número de gerações" configuration entry.
$ws.Range("A7").Value = "Número de gerações: 10000"

# Slight default-column-width recalculation that accompanied the edit
# in the authoring application.
$ws.StandardWidth = 12.00390625

# Move/leave the active selection on the newly added cell.
$ws.Range("A7").Select()
